$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing header labels (B1:D1) ---
$ws.Range("B1").Value = "32 Channels OS1 Normal, 45° FOV"
$ws.Range("C1").Value = "64 Channels OS1 Normal, 45° FOV"
$ws.Range("D1").Value = "128 Channels OS1 Normal, 45° FOV"

# --- Add new header labels (E1:H1) ---
$ws.Range("E1").Value = "32 Channels OS1 Below Horizon, 22.5° FOV"
$ws.Range("F1").Value = "64 Channels OS1 Below Horizon, 22.5° FOV"
$ws.Range("G1").Value = "128 Channels OS1 Below Horizon, 22.5° FOV"
$ws.Range("H1").Value = "128 Channels OS2 Normal, 22.5° FOV"

# Copy the header formatting (bold, border, centered) from D1 onto the new header cells
$ws.Range("D1").Copy()
$ws.Range("E1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update the data table with the new columns of cone-hit counts ---
$data = @{
    2 = @(0, 0, 0, 0, 0, 0, 0);
    3 = @(1, 3, 6, 3, 7, 13, 13);
    4 = @(1, 2, 4, 2, 4, 8, 9);
    5 = @(0, 1, 3, 1, 3, 5, 5);
    6 = @(1, 1, 1, 1, 1, 3, 3);
    7 = @(0, 1, 1, 0, 1, 2, 1);
    8 = @(0, 0, 1, 1, 1, 1, 2);
    9 = @(0, 0, 0, 0, 0, 0, 0)
}

$columns = @("B", "C", "D", "E", "F", "G", "H")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$row").Value = $values[$i]
    }
}
